$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '43.113.30'
$ws.Range("E2").Value = '  +5.15%  '
$ws.Range("D3").Value = '2.249.37'
$ws.Range("E3").Value = '  +4.39%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.00'
$ws.Range("E4").Value = '  -0.08%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '245.22'
$ws.Range("E5").Value = '  +4.29%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.615'
$ws.Range("E6").Value = '  +2.11%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '75.45'
$ws.Range("E7").Value = '  +9.67%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '1.00'
$ws.Range("E8").Value = '  -0.08%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.607'
$ws.Range("E9").Value = '  +8.17%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '40.78'
$ws.Range("E10").Value = '  +5.65%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0933'
$ws.Range("E11").Value = '  +3.61%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '6.97'
$ws.Range("E12").Value = '  +5.51%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.101'
$ws.Range("E13").Value = '  +1.43%  '
$ws.Range("D14").Value = '2.588.04'
$ws.Range("E14").Value = '  +4.06%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '14.67'
$ws.Range("E15").Value = '  +4.23%  '
$ws.Range("D16").Value = '2.249.11'
$ws.Range("E16").Value = '  +4.03%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.796'
$ws.Range("E17").Value = '  +2.43%  '
$ws.Range("D18").Value = '43.014.53'
$ws.Range("E18").Value = '  +5.40%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.0000105'
$ws.Range("E19").Value = '  +6.84%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '71.18'
$ws.Range("E20").Value = '  +2.57%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '5.98'
$ws.Range("E21").Value = '  +4.36%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '9.91'
$ws.Range("E22").Value = '  +7.86%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '230.45'
$ws.Range("E23").Value = '  +3.29%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '2.20'
$ws.Range("E24").Value = '  +17.73%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '1.00'
$ws.Range("E25").Value = '  +0.15%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '10.92'
$ws.Range("E26").Value = '  +3.68%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '3.43'
$ws.Range("E27").Value = '  -0.24%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '2.25'
$ws.Range("E28").Value = '  +3.64%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '38.86'
$ws.Range("E29").Value = '  +31.68%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '2.21'
$ws.Range("E30").Value = '  +2.21%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '173.10'
$ws.Range("E31").Value = '  +3.42%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '20.32'
$ws.Range("E32").Value = '  +3.51%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.0796'
$ws.Range("E33").Value = '  +7.01%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '5.31'
$ws.Range("E34").Value = '  +5.48%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.122'
$ws.Range("E35").Value = '  +2.75%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.109'
$ws.Range("E36").Value = '  +8.78%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '4.36'
$ws.Range("E37").Value = '  +8.57%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.0334'
$ws.Range("E38").Value = '  +20.94%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '13.19'
$ws.Range("E39").Value = '  +15.25%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '2.13'
$ws.Range("E40").Value = '  +5.16%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '5.50'
$ws.Range("E41").Value = '  +4.35%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.206'
$ws.Range("E42").Value = '  +11.24%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '59.71'
$ws.Range("E43").Value = '  +3.41%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '105.52'
$ws.Range("E44").Value = '  +10.00%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '8.73'
$ws.Range("E45").Value = '  +7.09%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.483'
$ws.Range("E46").Value = '  +32.21%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.0994'
$ws.Range("E47").Value = '  +4.43%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '2.40'
$ws.Range("E48").Value = '  +12.44%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '1.11'
$ws.Range("E49").Value = '  +4.82%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '1.15'
$ws.Range("E50").Value = '  +4.24%  '
$ws.Range("D51").Value = '2.460.53'
$ws.Range("E51").Value = '  +4.11%  '
